$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$r = $ws.Range("B16")
$r.HorizontalAlignment = -4108
$r.WrapText = $true
